$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.929.50'
$ws.Range('E2').Value = '  -2.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.859.93'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.80'
$ws.Range('E5').Value = '  -2.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5034'
$ws.Range('E7').Value = '  -3.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3709'
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07119'
$ws.Range('E9').Value = '  -1.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8819'
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.52'
$ws.Range('E11').Value = '  -2.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07564'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.857.15'
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.288'
$ws.Range('E14').Value = '  -2.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.03'
$ws.Range('E15').Value = '  -3.21%  '
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008379'
$ws.Range('E17').Value = '  -3.82%  '
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9995'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.984.09'
$ws.Range('E20').Value = '  -2.64%  '
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.108.93'
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('E23').Value = '  -3.36%  '
$ws.Range('E24').Value = '  -1.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.845'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '146.87'
$ws.Range('E26').Value = '  -4.14%  '
$ws.Range('E27').Value = '  -2.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.094'
$ws.Range('E28').Value = '  -3.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.44'
$ws.Range('E29').Value = '  -1.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.645'
$ws.Range('E30').Value = '  -3.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.661'
$ws.Range('E31').Value = '  -3.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09023'
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05114'
$ws.Range('E33').Value = '  -3.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.021'
$ws.Range('E34').Value = '  -4.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.145'
$ws.Range('E35').Value = '  -7.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7201'
$ws.Range('E36').Value = '  -7.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02034'
$ws.Range('E37').Value = '  -2.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.032'
$ws.Range('E38').Value = '  -0.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.450'
$ws.Range('E39').Value = '  -6.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.074'
$ws.Range('E40').Value = '  -1.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5266'
$ws.Range('E41').Value = '  -4.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.503'
$ws.Range('E42').Value = '  -2.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '115.02'
$ws.Range('E43').Value = '  +1.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.219'
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1464'
$ws.Range('E45').Value = '  -2.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9997'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4597'
$ws.Range('E47').Value = '  -3.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.931'
$ws.Range('E48').Value = '  -4.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.561'
$ws.Range('E49').Value = '  -3.14%  '
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('E51').Value = '  -4.00%  '
